# Add a 4th quiz/worksheet pair of columns (Q4, W4) to the gradebook,
# mirroring the existing Q1/W1 .. Q3/W3 columns (E:J), extending the
# used range from A1:J33 to A1:L33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("K1").Value = "Q4"
$ws.Range("L1").Value = "W4"

# Per-student Q4 (out of 4) / W4 (out of 50) values, rows 2-33
$kValues = @(4, 0, 4, 4, 2, 4, 2, 2, 4, 2, 4, 2, 2, 0, 2, 4, 4, 4, 4, 2, 4, 4, 2, 2, 2, 4, 2, 4, 2, 2, 4, 4)
$lValues = @(50, 0, 50, 50, 50, 50, 50, 50, 50, 42, 50, 50, 50, 50, 50, 45, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50, 50)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
}

# Match the author's final selection/active cell
$ws.Range("L33").Select()
